$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Row 17 ("repaymentstrategy") previously carried a stray "RBI (India)" value;
# update it to the new scenario value describing the repayment strategy.
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Scroll the view down and select the cell that was just updated, matching
# the author's recorded view state after making the change.
$ws.Activate()
$ws.Range("B17").Select()
$excel.ActiveWindow.ScrollRow = 10
